$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 701
$ws.Range("I18").Value = 701
$ws.Range("K18").Value = 701
$ws.Range("M18").Value = -417

$ws.Range("H138").Value = 3118.4119
$ws.Range("I138").Value = 3011.9333
$ws.Range("J138").Value = 3917
$ws.Range("K138").Value = 9035.7999
$ws.Range("L138").Value = 11751
$ws.Range("M138").Value = -3895.7999
$ws.Range("N138").Value = -22031

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 3572.1428
$ws.Range("I6").Value = 3666.6667
$ws.Range("J6").Value = 3501.25
$ws.Range("K6").Value = 3666.6667
$ws.Range("L6").Value = 3501.25
$ws.Range("M6").Value = -3493.6667
$ws.Range("N6").Value = -3847.25

$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

$ws.Range("H61").Value = 1754.7
$ws.Range("I61").Value = 2004.8572
$ws.Range("K61").Value = 2004.8572
$ws.Range("M61").Value = -1792.8572

$ws.Range("H110").Value = 1251.6666
$ws.Range("I110").Value = 1031
$ws.Range("J110").Value = 1693
$ws.Range("K110").Value = 1031
$ws.Range("L110").Value = 1693
$ws.Range("M110").Value = 1014
$ws.Range("N110").Value = -5783

$ws.Range("H136").Value = 1754.7
$ws.Range("I136").Value = 2004.8572
$ws.Range("K136").Value = 6014.571599999999
$ws.Range("M136").Value = -3464.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 9173.111000000001
$ws.Range("I7").Value = 551.5
$ws.Range("J7").Value = 16070.4
$ws.Range("K7").Value = 551.5
$ws.Range("L7").Value = 16070.4
$ws.Range("M7").Value = -438.5
$ws.Range("N7").Value = -16296.4

$ws.Range("H10").Value = 1077.1111
$ws.Range("J10").Value = 2304
$ws.Range("L10").Value = 2304
$ws.Range("N10").Value = -2584

$ws.Range("H16").Value = 1999
$ws.Range("J16").Value = 1999
$ws.Range("L16").Value = 1999
$ws.Range("N16").Value = -2339

$ws.Range("H20").Value = 1417.8
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H36").Value = 2749.6667
$ws.Range("I36").Value = 2749.6667
$ws.Range("K36").Value = 2749.6667
$ws.Range("M36").Value = -2215.6667

$ws.Range("H134").Value = 8000
$ws.Range("I134").Value = 8000
$ws.Range("K134").Value = 24000
$ws.Range("M134").Value = -21465

$ws.Range("H140").Value = 60000
$ws.Range("J140").Value = 60000
$ws.Range("L140").Value = 60000
$ws.Range("N140").Value = -70360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1247.5
$ws.Range("J3").Value = 1272.5
$ws.Range("L3").Value = 1272.5
$ws.Range("N3").Value = -1498.5

$ws.Range("H6").Value = 2275.75
$ws.Range("J6").Value = 3500
$ws.Range("L6").Value = 3500
$ws.Range("N6").Value = -3726

$ws.Range("H11").Value = 5495.778
$ws.Range("J11").Value = 6908.143
$ws.Range("L11").Value = 6908.143
$ws.Range("N11").Value = -7188.143

$ws.Range("H22").Value = 1600.2
$ws.Range("I22").Value = 999.5
$ws.Range("K22").Value = 999.5
$ws.Range("M22").Value = -649.5

$ws.Range("H31").Value = 1940.8334
$ws.Range("I31").Value = 1940.8334
$ws.Range("K31").Value = 1940.8334
$ws.Range("M31").Value = -1645.8334

$ws.Range("H34").Value = 1940.8334
$ws.Range("I34").Value = 1940.8334
$ws.Range("K34").Value = 1940.8334
$ws.Range("M34").Value = -1738.8334

$ws.Range("H35").Value = 8395.6
$ws.Range("I35").Value = 9237.5
$ws.Range("K35").Value = 9237.5
$ws.Range("M35").Value = -8943.5

$ws.Range("H105").Value = 397
$ws.Range("I105").Value = 399
$ws.Range("K105").Value = 399
$ws.Range("M105").Value = 1348

$ws.Range("H107").Value = 1021.8182
$ws.Range("I107").Value = 832.6667
$ws.Range("J107").Value = 1248.8
$ws.Range("K107").Value = 832.6667
$ws.Range("L107").Value = 1248.8
$ws.Range("M107").Value = 1087.3333
$ws.Range("N107").Value = -5088.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 819.2222
$ws.Range("I6").Value = 43.5
$ws.Range("J6").Value = 1439.8
$ws.Range("K6").Value = 130.5
$ws.Range("L6").Value = 4319.4
$ws.Range("M6").Value = -17.5
$ws.Range("N6").Value = -4545.4

$ws.Range("H26").Value = 1916
$ws.Range("J26").Value = 2374
$ws.Range("L26").Value = 7122
$ws.Range("N26").Value = -7698

$ws.Range("H46").Value = 5250.3335
$ws.Range("J46").Value = 9499.666999999999
$ws.Range("L46").Value = 28499.001
$ws.Range("N46").Value = -28681.001

$ws.Range("H55").Value = 3173.8
$ws.Range("J55").Value = 3912.5
$ws.Range("L55").Value = 11737.5
$ws.Range("N55").Value = -12091.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 1000
$ws.Range("J23").Value = 1000
$ws.Range("L23").Value = 1000
$ws.Range("N23").Value = -1446

$ws.Range("H59").Value = 45000
$ws.Range("J59").Value = 45000
$ws.Range("L59").Value = 45000
$ws.Range("N59").Value = -46166

$ws.Range("H95").Value = 10000
$ws.Range("J95").Value = 10000
$ws.Range("L95").Value = 10000
$ws.Range("N95").Value = -15492

$ws.Range("H113").Value = 892.5
$ws.Range("I113").Value = 785
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 785
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1385
$ws.Range("N113").Value = -5340

$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 10000
$ws.Range("K126").Value = 30000
$ws.Range("M126").Value = -27530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("N30").ClearContents()

$ws.Range("H35").Value = 15249.375
$ws.Range("I35").Value = 7000
$ws.Range("J35").Value = 20199
$ws.Range("K35").Value = 7000
$ws.Range("L35").Value = 20199
$ws.Range("M35").Value = -6664
$ws.Range("N35").Value = -20871

$ws.Range("H132").Value = 2924.7778
$ws.Range("I132").Value = 2627.5715
$ws.Range("K132").Value = 7882.7145
$ws.Range("M132").Value = -5352.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 10000
$ws.Range("J2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("N2").Value = -10224

$ws.Range("H4").Value = 1989.9
$ws.Range("J4").Value = 1612.5
$ws.Range("L4").Value = 1612.5
$ws.Range("N4").Value = -1838.5

$ws.Range("H5").Value = 2502000
$ws.Range("I5").Value = 5075000
$ws.Range("J5").Value = 786666.7
$ws.Range("K5").Value = 5075000
$ws.Range("L5").Value = 786666.7
$ws.Range("M5").Value = -5074888
$ws.Range("N5").Value = -786890.7

$ws.Range("H6").Value = 14289042
$ws.Range("J6").Value = 947
$ws.Range("L6").Value = 947
$ws.Range("N6").Value = -1177

$ws.Range("H10").Value = 10005
$ws.Range("I10").Value = 10005
$ws.Range("K10").Value = 10005
$ws.Range("M10").Value = -9836

$ws.Range("H26").Value = 18998
$ws.Range("J26").Value = 18998
$ws.Range("L26").Value = 18998
$ws.Range("N26").Value = -19584

$ws.Range("H107").Value = 1091.3636
$ws.Range("I107").Value = 888.8889
$ws.Range("K107").Value = 2666.6667
$ws.Range("M107").Value = -746.6667000000002

$ws.Range("H113").Value = 416.44446
$ws.Range("I113").Value = 389.66666
$ws.Range("J113").Value = 470
$ws.Range("K113").Value = 1168.99998
$ws.Range("L113").Value = 1410
$ws.Range("M113").Value = 1001.00002
$ws.Range("N113").Value = -5750

$ws.Range("H136").Value = 2541.2856
$ws.Range("I136").Value = 2541.2856
$ws.Range("K136").Value = 7623.8568
$ws.Range("M136").Value = -5073.8568
